# "corrected hero attack path"
# The hero-attack / defence section of the battle algorithm diagram is
# reshuffled & extended: a missing "hero defence check" step is inserted,
# "defense" is renamed to "defence" throughout, and a whole new
# "monster defence" step (4) is appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# ---------------------------------------------------------------------
# 1) Shift the whole "hero received dmg / armor dmg / game over / hero
#    attacks..." block down to make room for the new "hero deffence
#    check" row, working bottom-up so sources are never clobbered
#    before they are read. Copy() preserves both value and style, and
#    only reuses already-existing shared strings (it does not change
#    shared-string allocation order).
# ---------------------------------------------------------------------
$ws.Range("B35:C35").Copy($ws.Range("B37"))
$ws.Range("B33:C33").Copy($ws.Range("B35"))
$ws.Range("B31:C31").Copy($ws.Range("B33"))
$ws.Range("B29:C29").Copy($ws.Range("B31"))
$ws.Range("B27:C27").Copy($ws.Range("B29"))
$ws.Range("C25").Copy($ws.Range("C26"))
$ws.Range("B23").Copy($ws.Range("C24"))

# Clear the now-vacated source cells that were not overwritten by a
# move (Clear, not ClearContents, so the now-empty row 27 disappears
# entirely instead of leaving a blank styled row behind).
$ws.Range("B23").Clear()
$ws.Range("C25").Clear()
$ws.Range("B27:C27").Clear()

# ---------------------------------------------------------------------
# 2) Fill in all brand-new text, in first-use order so the regenerated
#    shared-strings table lines up with the authored workbook.
# ---------------------------------------------------------------------
$ws.Range("B23").Value = 'if "hero deffence" < "actual monster attack"'
$ws.Range("B39").Value = "monster deffence"
$ws.Range("C39").Value = "monster.deffence"
$ws.Range("C43").Value = "monster received dmg = actual hero attack - monster deffence"
$ws.Range("B37").Value = "4)"

# Give the new "hero chooses defence" label the same formatting as its
# "hero chooses attack" sibling by copying it first, then overwrite
# the copied text with the real label.
$ws.Range("B11").Copy($ws.Range("M11"))
$ws.Range("M11").Value = "hero chooses defence"

$ws.Range("B41").Value = 'if "actual hero attack" < "monster defence"'
$ws.Range("B6").Value = "defence;"
$ws.Range("F6").Value = "defence;"
$ws.Range("C21").Value = "hero.defence + hero.armor.head + hero.armor.shield + hero.armor.chestPlate + hero.armor.shoes"

# ---------------------------------------------------------------------
# 3) Update the view to match what the author left it at
# ---------------------------------------------------------------------
$ws.Range("C39").Select()
$excel.ActiveWindow.ScrollRow = 19
